# Auto-generated edit script: apply numeric corrections to profit-tracking sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 92891.17999999999
$ws.Cells.Item(28, 9).Value = 112210.78
$ws.Cells.Item(28, 10).Value = 5953
$ws.Cells.Item(28, 11).Value = 112210.78
$ws.Cells.Item(28, 12).Value = 5953
$ws.Cells.Item(28, 13).Value = -111725.78
$ws.Cells.Item(28, 14).Value = -6923

$ws.Cells.Item(132, 8).Value = 1030.5869
$ws.Cells.Item(132, 9).Value = 967.55817
$ws.Cells.Item(132, 11).Value = 2902.67451
$ws.Cells.Item(132, 13).Value = -372.6745099999998

$ws.Cells.Item(135, 8).Value = 789.5833
$ws.Cells.Item(135, 9).Value = 651.7778
$ws.Cells.Item(135, 10).Value = 1203
$ws.Cells.Item(135, 11).Value = 5866.000199999999
$ws.Cells.Item(135, 12).Value = 10827
$ws.Cells.Item(135, 13).Value = -3331.000199999999
$ws.Cells.Item(135, 14).Value = -15897

$ws.Cells.Item(138, 8).Value = 2849.6428
$ws.Cells.Item(138, 9).Value = 2323.5557
$ws.Cells.Item(138, 10).Value = 3244.2083
$ws.Cells.Item(138, 11).Value = 6970.6671
$ws.Cells.Item(138, 12).Value = 9732.624899999999
$ws.Cells.Item(138, 13).Value = -1830.6671
$ws.Cells.Item(138, 14).Value = -20012.6249

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2920.32
$ws.Cells.Item(132, 9).Value = 2325.2632
$ws.Cells.Item(132, 10).Value = 4804.6665
$ws.Cells.Item(132, 11).Value = 6975.7896
$ws.Cells.Item(132, 12).Value = 14413.9995
$ws.Cells.Item(132, 13).Value = -4445.7896
$ws.Cells.Item(132, 14).Value = -19473.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2952
$ws.Cells.Item(20, 9).Value = 2226.3
$ws.Cells.Item(20, 10).Value = 4403.4
$ws.Cells.Item(20, 11).Value = 2226.3
$ws.Cells.Item(20, 12).Value = 4403.4
$ws.Cells.Item(20, 13).Value = -1979.3
$ws.Cells.Item(20, 14).Value = -4897.4

$ws.Cells.Item(140, 8).Value = 69999
$ws.Cells.Item(140, 10).Value = 69999
$ws.Cells.Item(140, 12).Value = 69999
$ws.Cells.Item(140, 14).Value = -80359

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 4757.0625
$ws.Cells.Item(122, 10).Value = 7821.75
$ws.Cells.Item(122, 12).Value = 23465.25
$ws.Cells.Item(122, 14).Value = -28365.25

$ws.Cells.Item(134, 8).Value = 3435.5
$ws.Cells.Item(134, 9).Value = 2227.2727
$ws.Cells.Item(134, 10).Value = 6093.6
$ws.Cells.Item(134, 11).Value = 6681.8181
$ws.Cells.Item(134, 12).Value = 18280.8
$ws.Cells.Item(134, 13).Value = -4146.8181
$ws.Cells.Item(134, 14).Value = -23350.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2127.95
$ws.Cells.Item(5, 9).Value = 787.7273
$ws.Cells.Item(5, 10).Value = 3766
$ws.Cells.Item(5, 11).Value = 2363.1819
$ws.Cells.Item(5, 12).Value = 11298
$ws.Cells.Item(5, 13).Value = -2251.1819
$ws.Cells.Item(5, 14).Value = -11522

$ws.Cells.Item(9, 8).Value = 354701.88
$ws.Cells.Item(9, 9).Value = 577141.2
$ws.Cells.Item(9, 11).Value = 1731423.6
$ws.Cells.Item(9, 13).Value = -1731199.6

$ws.Cells.Item(12, 8).Value = 80.2
$ws.Cells.Item(12, 9).Value = 100.666664
$ws.Cells.Item(12, 10).Value = 71.42856999999999
$ws.Cells.Item(12, 11).Value = 301.999992
$ws.Cells.Item(12, 12).Value = 214.28571
$ws.Cells.Item(12, 13).Value = -128.999992
$ws.Cells.Item(12, 14).Value = -560.28571

$ws.Cells.Item(15, 8).Value = 50.75
$ws.Cells.Item(15, 9).Value = 50.75
$ws.Cells.Item(15, 11).Value = 152.25
$ws.Cells.Item(15, 13).Value = -12.25

$ws.Cells.Item(16, 8).Value = 1191.6666
$ws.Cells.Item(16, 10).Value = 1633.3334
$ws.Cells.Item(16, 12).Value = 4900.0002
$ws.Cells.Item(16, 14).Value = -5246.0002

$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).ClearContents()
$ws.Cells.Item(20, 14).ClearContents()

$ws.Cells.Item(21, 8).Value = 100.666664
$ws.Cells.Item(21, 9).Value = 100.666664
$ws.Cells.Item(21, 11).Value = 301.999992
$ws.Cells.Item(21, 13).Value = -128.999992

$ws.Cells.Item(22, 8).Value = 7185.5713
$ws.Cells.Item(22, 9).Value = 300
$ws.Cells.Item(22, 11).Value = 900
$ws.Cells.Item(22, 13).Value = -731

$ws.Cells.Item(27, 8).Value = 7185.5713
$ws.Cells.Item(27, 9).Value = 300
$ws.Cells.Item(27, 11).Value = 900
$ws.Cells.Item(27, 13).Value = -798

$ws.Cells.Item(32, 8).Value = 9118.666999999999
$ws.Cells.Item(32, 10).Value = 9118.666999999999
$ws.Cells.Item(32, 12).Value = 27356.001
$ws.Cells.Item(32, 14).Value = -27922.001

$ws.Cells.Item(39, 8).Value = 3522.4443
$ws.Cells.Item(39, 9).Value = 3500
$ws.Cells.Item(39, 10).Value = 3525.25
$ws.Cells.Item(39, 11).Value = 10500
$ws.Cells.Item(39, 12).Value = 10575.75
$ws.Cells.Item(39, 13).Value = -10206
$ws.Cells.Item(39, 14).Value = -11163.75

$ws.Cells.Item(40, 8).Value = 261.08334
$ws.Cells.Item(40, 9).Value = 27.25
$ws.Cells.Item(40, 10).Value = 378
$ws.Cells.Item(40, 11).Value = 109
$ws.Cells.Item(40, 12).Value = 1512
$ws.Cells.Item(40, 13).Value = -40
$ws.Cells.Item(40, 14).Value = -1650

$ws.Cells.Item(44, 8).Value = 400
$ws.Cells.Item(44, 9).Value = 400
$ws.Cells.Item(44, 11).Value = 1200
$ws.Cells.Item(44, 13).Value = -802

$ws.Cells.Item(46, 8).Value = 1873.25
$ws.Cells.Item(46, 9).Value = 899.5
$ws.Cells.Item(46, 10).Value = 2847
$ws.Cells.Item(46, 11).Value = 2698.5
$ws.Cells.Item(46, 12).Value = 8541
$ws.Cells.Item(46, 13).Value = -2607.5
$ws.Cells.Item(46, 14).Value = -8723

$ws.Cells.Item(135, 8).Value = 2127.95
$ws.Cells.Item(135, 9).Value = 787.7273
$ws.Cells.Item(135, 10).Value = 3766
$ws.Cells.Item(135, 11).Value = 7089.545700000001
$ws.Cells.Item(135, 12).Value = 33894
$ws.Cells.Item(135, 13).Value = -4554.545700000001
$ws.Cells.Item(135, 14).Value = -38964

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()

$ws.Cells.Item(70, 8).Value = 15385.728
$ws.Cells.Item(70, 9).Value = 5593.1875
$ws.Cells.Item(70, 11).Value = 5593.1875
$ws.Cells.Item(70, 13).Value = -5323.1875

$ws.Cells.Item(73, 8).Value = 15385.728
$ws.Cells.Item(73, 9).Value = 5593.1875
$ws.Cells.Item(73, 11).Value = 5593.1875
$ws.Cells.Item(73, 13).Value = -4657.1875

$ws.Cells.Item(133, 8).Value = 69995.42999999999
$ws.Cells.Item(133, 10).Value = 69995.42999999999
$ws.Cells.Item(133, 12).Value = 69995.42999999999
$ws.Cells.Item(133, 14).Value = -80115.42999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2627.25
$ws.Cells.Item(93, 9).Value = 2209.95
$ws.Cells.Item(93, 11).Value = 2209.95
$ws.Cells.Item(93, 13).Value = -961.9499999999998

$ws.Cells.Item(132, 8).Value = 2746.6667
$ws.Cells.Item(132, 9).Value = 1235.9524
$ws.Cells.Item(132, 10).Value = 4861.6665
$ws.Cells.Item(132, 11).Value = 3707.857199999999
$ws.Cells.Item(132, 12).Value = 14584.9995
$ws.Cells.Item(132, 13).Value = -1177.857199999999
$ws.Cells.Item(132, 14).Value = -19644.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 7681.9287
$ws.Cells.Item(122, 9).Value = 2126.4666
$ws.Cells.Item(122, 10).Value = 14092.077
$ws.Cells.Item(122, 11).Value = 6379.399800000001
$ws.Cells.Item(122, 12).Value = 42276.231
$ws.Cells.Item(122, 13).Value = -3929.399800000001
$ws.Cells.Item(122, 14).Value = -47176.231

$ws.Cells.Item(132, 8).Value = 5898.68
$ws.Cells.Item(132, 9).Value = 5436.15
$ws.Cells.Item(132, 10).Value = 7748.8
$ws.Cells.Item(132, 11).Value = 16308.45
$ws.Cells.Item(132, 12).Value = 23246.4
$ws.Cells.Item(132, 13).Value = -13778.45
$ws.Cells.Item(132, 14).Value = -28306.4

$ws.Cells.Item(136, 8).Value = 4521.6665
$ws.Cells.Item(136, 9).Value = 2216.923
$ws.Cells.Item(136, 11).Value = 6650.768999999999
$ws.Cells.Item(136, 13).Value = -4100.768999999999
